$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.164.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.38%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("D9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

# Row 11
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "

# Row 12
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "

# Row 13
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.627.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.77%  "

# Row 14
$ws.Range("E14").Value = "  +0.50%  "

# Row 15
$ws.Range("E15").Value = "  +2.19%  "

# Row 16
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.178.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.19%  "

# Row 18
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "

# Row 19
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.33%  "

# Row 21
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.43%  "

# Row 22
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "

# Row 25
$ws.Range("D25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27
$ws.Range("E27").Value = "  +1.02%  "

# Row 28
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("E29").Value = "  -0.74%  "

# Row 30
$ws.Range("E30").Value = "  +0.81%  "

# Row 31
$ws.Range("E31").Value = "  +0.43%  "

# Row 32
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "

# Row 33
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.305.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.90%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.75%  "

# Row 36
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "

# Row 37
$ws.Range("E37").Value = "  -1.34%  "

# Row 38
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.84%  "

# Row 39
$ws.Range("D39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.85%  "

# Row 40
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("E42").Value = "  +5.96%  "

# Row 43
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.780.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("D46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "

# Row 47
$ws.Range("D47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.47%  "

# Row 48
$ws.Range("D48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "

# Row 49
$ws.Range("D49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "

# Row 50
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "

# Row 51
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0963"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
